$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.532141
$ws.Range("H2").Value = 4.596423
$ws.Range("I2").Value = 0.08900664250669833
$ws.Range("J2").Value = 0.08900664250669831
$ws.Range("Q2").Value = 0.03901801341966666
$ws.Range("R2").Value = 0.3511621207769999
$ws.Range("S2").Value = 0.08900664250669833
$ws.Range("T2").Value = 0.08900664250669831

# Row 3 (FAPs)
$ws.Range("I3").Value = 0.6169137955113024
$ws.Range("J3").Value = 0.6169137955113023
$ws.Range("S3").Value = 0.6169137955113024
$ws.Range("T3").Value = 0.6169137955113023

# Row 4 (MuSCs)
$ws.Range("G4").Value = 4.902263666666666
$ws.Range("H4").Value = 14.706791
$ws.Range("I4").Value = 0.2847871244569372
$ws.Range("J4").Value = 0.2847871244569371
$ws.Range("Q4").Value = 0.1248426806232222
$ws.Range("R4").Value = 1.123584125609
$ws.Range("S4").Value = 0.2847871244569372
$ws.Range("T4").Value = 0.2847871244569371

# Row 5 (Resolving-Mac)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.159958
$ws.Range("H5").Value = 0.479874
$ws.Range("I5").Value = 0.009292437525062282
$ws.Range("J5").Value = 0.009292437525062281
$ws.Range("Q5").Value = 0.004073543747333333
$ws.Range("R5").Value = 0.036661893726
$ws.Range("S5").Value = 0.009292437525062282
$ws.Range("T5").Value = 0.009292437525062281
